$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 15500
$ws.Range("I7").Value = 12000
$ws.Range("J7").Value = 19000
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 19000
$ws.Range("M7").Value = -11888
$ws.Range("N7").Value = -19224

$ws.Range("H14").Value = 15500
$ws.Range("I14").Value = 12000
$ws.Range("J14").Value = 19000
$ws.Range("K14").Value = 12000
$ws.Range("L14").Value = 19000
$ws.Range("M14").Value = -11809
$ws.Range("N14").Value = -19382

$ws.Range("H33").Value = 2525.9285
$ws.Range("J33").Value = 6236.8
$ws.Range("L33").Value = 6236.8
$ws.Range("N33").Value = -6694.8

$ws.Range("H98").Value = 1174.4736
$ws.Range("I98").Value = 749.3333
$ws.Range("J98").Value = 2768.75
$ws.Range("K98").Value = 749.3333
$ws.Range("L98").Value = 2768.75
$ws.Range("M98").Value = 748.6667
$ws.Range("N98").Value = -5764.75

$ws.Range("H122").Value = 1174.4736
$ws.Range("I122").Value = 749.3333
$ws.Range("J122").Value = 2768.75
$ws.Range("K122").Value = 2247.9999
$ws.Range("L122").Value = 8306.25
$ws.Range("M122").Value = 202.0001000000002
$ws.Range("N122").Value = -13206.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5853.4414
$ws.Range("I32").Value = 2724.6948
$ws.Range("J32").Value = 26364.111
$ws.Range("K32").Value = 2724.6948
$ws.Range("L32").Value = 26364.111
$ws.Range("M32").Value = -2437.6948
$ws.Range("N32").Value = -26938.111

$ws.Range("H61").Value = 5425.7095
$ws.Range("I61").Value = 5454.9546
$ws.Range("K61").Value = 5454.9546
$ws.Range("M61").Value = -5242.9546

$ws.Range("H132").Value = 4744.7915
$ws.Range("I132").Value = 2651.7058
$ws.Range("K132").Value = 7955.117400000001
$ws.Range("M132").Value = -5425.117400000001

$ws.Range("H136").Value = 5425.7095
$ws.Range("I136").Value = 5454.9546
$ws.Range("K136").Value = 16364.8638
$ws.Range("M136").Value = -13814.8638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1418.5862
$ws.Range("I86").Value = 1219.381
$ws.Range("J86").Value = 1941.5
$ws.Range("K86").Value = 1219.381
$ws.Range("L86").Value = 1941.5
$ws.Range("M86").Value = -96.38100000000009
$ws.Range("N86").Value = -4187.5

$ws.Range("H89").Value = 1418.5862
$ws.Range("I89").Value = 1219.381
$ws.Range("J89").Value = 1941.5
$ws.Range("K89").Value = 6096.905000000001
$ws.Range("L89").Value = 9707.5
$ws.Range("M89").Value = -480.9050000000007
$ws.Range("N89").Value = -20939.5

$ws.Range("H107").Value = 1990.8334
$ws.Range("I107").Value = 1750.625
$ws.Range("K107").Value = 1750.625
$ws.Range("M107").Value = 169.375

$ws.Range("H134").Value = 4348.1904
$ws.Range("J134").Value = 11828.5
$ws.Range("L134").Value = 35485.5
$ws.Range("N134").Value = -40555.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 199500
$ws.Range("J63").Value = 199500
$ws.Range("L63").Value = 199500
$ws.Range("N63").Value = -200872

$ws.Range("H66").Value = 199500
$ws.Range("J66").Value = 199500
$ws.Range("L66").Value = 598500
$ws.Range("N66").Value = -605364

$ws.Range("H75").Value = 199500
$ws.Range("J75").Value = 199500
$ws.Range("L75").Value = 199500
$ws.Range("N75").Value = -201496

$ws.Range("H78").Value = 199500
$ws.Range("J78").Value = 199500
$ws.Range("L78").Value = 598500
$ws.Range("N78").Value = -608484

$ws.Range("H81").Value = 199500
$ws.Range("J81").Value = 199500
$ws.Range("L81").Value = 199500
$ws.Range("N81").Value = -201496

$ws.Range("H84").Value = 199500
$ws.Range("J84").Value = 199500
$ws.Range("L84").Value = 598500
$ws.Range("N84").Value = -608484

$ws.Range("H87").Value = 199500
$ws.Range("J87").Value = 199500
$ws.Range("L87").Value = 199500
$ws.Range("N87").Value = -201872

$ws.Range("H90").Value = 199500
$ws.Range("J90").Value = 199500
$ws.Range("L90").Value = 598500
$ws.Range("N90").Value = -610356

$ws.Range("H132").Value = 3654.375
$ws.Range("J132").Value = 7160
$ws.Range("L132").Value = 21480
$ws.Range("N132").Value = -26540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5.3333335
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H112").Value = 19542.5
$ws.Range("I112").Value = 17465
$ws.Range("J112").Value = 19958
$ws.Range("K112").Value = 52395
$ws.Range("L112").Value = 59874
$ws.Range("M112").Value = -51287
$ws.Range("N112").Value = -62090

$ws.Range("H137").Value = 3207.5557
$ws.Range("I137").Value = 2113.7693
$ws.Range("J137").Value = 3825.7827
$ws.Range("K137").Value = 6341.3079
$ws.Range("L137").Value = 11477.3481
$ws.Range("M137").Value = -1241.3079
$ws.Range("N137").Value = -21677.3481

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 32200
$ws.Range("I14").Value = 32200
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 32200
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -32032

$ws.Range("H122").Value = 4660.448
$ws.Range("I122").Value = 3091.6667
$ws.Range("K122").Value = 9275.000100000001
$ws.Range("M122").Value = -6825.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 7503.5
$ws.Range("J13").Value = 7503.5
$ws.Range("L13").Value = 7503.5
$ws.Range("N13").Value = -7783.5

$ws.Range("H40").Value = 6118.45
$ws.Range("I40").Value = 5222.933
$ws.Range("J40").Value = 8805
$ws.Range("K40").Value = 5222.933
$ws.Range("L40").Value = 8805
$ws.Range("M40").Value = -5086.933
$ws.Range("N40").Value = -9077

$ws.Range("H132").Value = 5401.1816
$ws.Range("I132").Value = 4228.2104
$ws.Range("K132").Value = 12684.6312
$ws.Range("M132").Value = -10154.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3771.8235
$ws.Range("I132").Value = 2741.4
$ws.Range("J132").Value = 11500
$ws.Range("K132").Value = 8224.200000000001
$ws.Range("L132").Value = 34500
$ws.Range("M132").Value = -5694.200000000001
$ws.Range("N132").Value = -39560
